# Applies a cyclic rotation of the data in rows 2-4 (row2 -> row3 -> row4 -> row2)
# i.e. the new row2 gets the old row4's data, new row3 gets the old row2's data,
# and new row4 gets the old row3's data. Columns that are identical across all
# three rows (C, K, N, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AF, AG, AT, AW, AX, AY)
# are left untouched since the edit is invisible there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 2 (was old row 4: Vätteros / Lathraea squamaria) ---
$ws.Range("A2").Value = 92877460
$ws.Range("B2").Value = 104490
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 219686
$ws.Range("F2").Value = "Vätteros"
$ws.Range("G2").Value = "Lathraea squamaria"
$ws.Range("H2").Value = "L."
$ws.Range("I2").Value = "'20"
$ws.Range("J2").Value = "plantor/tuvor"
$ws.Range("Q2").Value = 563580.4029258011
$ws.Range("R2").Value = 6434562.017619756

# --- New row 3 (was old row 2: Tallticka / Porodaedalea pini) ---
$ws.Range("A3").Value = 92877504
$ws.Range("B3").Value = 89412
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 5442
$ws.Range("F3").Value = "Tallticka"
$ws.Range("G3").Value = "Porodaedalea pini"
$ws.Range("H3").Value = "(Brot.) Murrill"
$ws.Range("AC3").Value = ""
$ws.Range("Q3").Value = 563565.9699459416
$ws.Range("R3").Value = 6434535.848147285

# --- New row 4 (was old row 3: Vedskivlav / Hertelidea botryosa) ---
$ws.Range("A4").Value = 92877498
$ws.Range("B4").Value = 78098
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6453
$ws.Range("F4").Value = "Vedskivlav"
$ws.Range("G4").Value = "Hertelidea botryosa"
$ws.Range("H4").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("AC4").Value = "murken tallstubbe. Artbestämd av Steve Daurer."
$ws.Range("Q4").Value = 563565.9699459416
$ws.Range("R4").Value = 6434535.848147285
